# Update cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue 'D2' '30.206.53'
$ws.Range('E2').Value = '  +6.21%  '
Set-TextValue 'D3' '1.918.87'
$ws.Range('E3').Value = '  +2.86%  '
Set-TextValue 'D4' '1.001'
$ws.Range('E4').Value = '  -0.58%  '
Set-TextValue 'D5' '330.40'
$ws.Range('E5').Value = '  +5.14%  '
Set-TextValue 'D6' '0.9999'
$ws.Range('E6').Value = '  -0.65%  '
Set-TextValue 'D7' '0.5218'
$ws.Range('E7').Value = '  +3.03%  '
Set-TextValue 'D8' '0.4093'
$ws.Range('E8').Value = '  +5.07%  '
Set-TextValue 'D9' '0.08528'
$ws.Range('E9').Value = '  +2.69%  '
Set-TextValue 'B10' 'OKB'
Set-TextValue 'C10' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D10' '42.96'
$ws.Range('E10').Value = '  +1.40%  '
Set-TextValue 'B11' 'Polygon'
Set-TextValue 'C11' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D11' '1.128'
$ws.Range('E11').Value = '  +2.27%  '
Set-TextValue 'D12' '23.14'
$ws.Range('E12').Value = '  +14.20%  '
Set-TextValue 'D13' '6.454'
$ws.Range('E13').Value = '  +4.54%  '
Set-TextValue 'D14' '1.914.93'
$ws.Range('E14').Value = '  +2.82%  '
Set-TextValue 'D15' '7.406'
$ws.Range('E15').Value = '  +2.53%  '
Set-TextValue 'D16' '1.001'
$ws.Range('E16').Value = '  -0.59%  '
Set-TextValue 'D17' '95.27'
$ws.Range('E17').Value = '  +4.71%  '
Set-TextValue 'D18' '0.00001115'
$ws.Range('E18').Value = '  +1.68%  '
Set-TextValue 'D19' '0.06697'
$ws.Range('E19').Value = '  -0.24%  '
Set-TextValue 'D20' '18.51'
$ws.Range('E20').Value = '  +5.48%  '
Set-TextValue 'D21' '0.9994'
$ws.Range('E21').Value = '  -0.66%  '
Set-TextValue 'D22' '6.012'
$ws.Range('E22').Value = '  +2.04%  '
Set-TextValue 'D23' '30.218.77'
$ws.Range('E23').Value = '  +6.13%  '
$ws.Range('E24').Value = '  +2.92%  '
Set-TextValue 'D25' '2.224'
$ws.Range('E25').Value = '  +1.43%  '
Set-TextValue 'D26' '2.133.18'
$ws.Range('E26').Value = '  +2.85%  '
Set-TextValue 'D27' '161.27'
$ws.Range('E27').Value = '  +2.03%  '
Set-TextValue 'D28' '21.16'
$ws.Range('E28').Value = '  +2.93%  '
Set-TextValue 'D29' '2.422'
$ws.Range('E29').Value = '  +0.65%  '
Set-TextValue 'D30' '129.10'
$ws.Range('E30').Value = '  +2.35%  '
Set-TextValue 'D31' '1.087'
$ws.Range('E31').Value = '  +5.12%  '
Set-TextValue 'D32' '0.1070'
$ws.Range('E32').Value = '  +3.51%  '
Set-TextValue 'D33' '6.023'
$ws.Range('E33').Value = '  +4.61%  '
Set-TextValue 'D34' '3.613'
$ws.Range('E34').Value = '  -0.18%  '
Set-TextValue 'D35' '0.02493'
$ws.Range('E35').Value = '  +2.30%  '
Set-TextValue 'D36' '0.06584'
$ws.Range('E36').Value = '  +0.19%  '
Set-TextValue 'D37' '0.2213'
$ws.Range('E37').Value = '  +2.90%  '
Set-TextValue 'D38' '1.231'
$ws.Range('E38').Value = '  +4.46%  '
Set-TextValue 'D39' '5.180'
$ws.Range('E39').Value = '  +3.40%  '
Set-TextValue 'D40' '8.816'
$ws.Range('E40').Value = '  -1.36%  '
Set-TextValue 'B41' 'Aptos'
Set-TextValue 'C41' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D41' '11.81'
$ws.Range('E41').Value = '  +6.68%  '
Set-TextValue 'B42' 'TheSandbox'
Set-TextValue 'C42' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D42' '0.6536'
$ws.Range('E42').Value = '  +3.12%  '
Set-TextValue 'D43' '1.242'
$ws.Range('E43').Value = '  +0.73%  '
Set-TextValue 'D44' '0.6169'
$ws.Range('E44').Value = '  +3.38%  '
Set-TextValue 'D45' '13.34'
$ws.Range('E45').Value = '  +2.30%  '
Set-TextValue 'D46' '3.746'
$ws.Range('E46').Value = '  +1.93%  '
Set-TextValue 'D47' '2.089'
$ws.Range('E47').Value = '  +4.96%  '
Set-TextValue 'D48' '1.244'
$ws.Range('E48').Value = '  +3.36%  '
Set-TextValue 'D49' '124.30'
$ws.Range('E49').Value = '  +2.00%  '
Set-TextValue 'D50' '1.168'
$ws.Range('E50').Value = '  +1.49%  '
Set-TextValue 'D51' '79.73'
$ws.Range('E51').Value = '  +5.38%  '
